$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at 76-78; existing rows 76-95 shift down to 79-98.
$ws.Rows("76:78").Insert()

# Common columns (same across the 3 new "Carson" rows).
$ws.Range("A76:A78").Value = 2
$ws.Range("B76:B78").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C76:C78").Value = "Coquimbo"
$ws.Range("D76:D78").Value = 44609
$ws.Range("E76:E78").Value = 4
$ws.Range("F76:F78").Value = "Fruta"
$ws.Range("G76:G78").Value = 100103
$ws.Range("H76:H78").Value = "Frutos de hueso (carozo)"
$ws.Range("I76:I78").Value = 100103004
$ws.Range("J76:J78").Value = "Durazno"
$ws.Range("K76:K78").Value = "Carson"
$ws.Range("Q76:Q78").Value = "$/caja 16 kilos empedrada"
$ws.Range("R76:R78").Value = "Región de O'Higgins"
$ws.Range("T76:T78").Value = 16

# Row 76 - Especial
$ws.Range("L76").Value = "Especial"
$ws.Range("M76").Value = 200
$ws.Range("N76").Value = 17000
$ws.Range("O76").Value = 18000
$ws.Range("P76").Value = 17500
$ws.Range("S76").Value = 1094

# Row 77 - Primera
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 200
$ws.Range("N77").Value = 15000
$ws.Range("O77").Value = 16000
$ws.Range("P77").Value = 15500
$ws.Range("S77").Value = 969

# Row 78 - Segunda
$ws.Range("L78").Value = "Segunda"
$ws.Range("M78").Value = 200
$ws.Range("N78").Value = 13000
$ws.Range("O78").Value = 14000
$ws.Range("P78").Value = 13500
$ws.Range("S78").Value = 844
